$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODO CMS")

# New row 10: new TODO entry, status "offen", assignee "Jonas"
$ws.Range("A10").Value = "Bei Veränderung des Amounts von 0 auf >0 -> Mail an alle Vorbesteller senden, dass Produkt jetzt verfügbar und bestellt ist"
$ws.Range("B10").Value = "offen"
$ws.Range("B10").Style = "Schlecht"
$ws.Range("C10").Value = "Jonas"

# Update the active selection to B12 to match the recorded view state
$ws.Activate()
$ws.Range("B12").Select()
